# "New prompt format with order" - insert the missing "doze" word-family
# rows (doze/dozed/dozes/dozing) into the Sleep keyword dictionary so the
# list stays alphabetically ordered ahead of "estivate".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "estivate" currently starts at row 490; push it (and everything below)
# down by 4 rows to make room for the new "doze" entries.
$ws.Range("A490:A493").EntireRow.Insert()

$newRows = @(
    @("doze",   "Sleep", "doze"),
    @("dozed",  "Sleep", "doze"),
    @("dozes",  "Sleep", "doze"),
    @("dozing", "Sleep", "doze")
)

$r = 490
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r = $r + 1
}
